$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"

$ws.Range("G2").Value = "iaest-measure:corine-land-cover-2000-nivel-1-descripcion"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

$ws.Range("G5").Clear()
